$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 10
$ws.Range("C5").Value = 25

$ws.Range("C3").Select()
